# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 20908
$ws1.Range("F4").Value = 336
$ws1.Range("F7").Value = 7811
$ws1.Range("F8").Value = 545
$ws1.Range("F13").Value = 182
$ws1.Range("F19").Value = 503
$ws1.Range("F26").Value = 1163
$ws1.Range("F29").Value = 209
$ws1.Range("F32").Value = 118
$ws1.Range("F33").Value = 4978
$ws1.Range("F38").Value = 12952
$ws1.Range("F39").Value = 1353
$ws1.Range("F40").Value = 115
$ws1.Range("F43").Value = 296
$ws1.Range("F44").Value = 410

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 20908
$ws4.Range("F4").Value = 336
$ws4.Range("F7").Value = 7811
$ws4.Range("F8").Value = 545
$ws4.Range("F13").Value = 182
$ws4.Range("F19").Value = 503
$ws4.Range("F26").Value = 1163
$ws4.Range("F29").Value = 209
$ws4.Range("F33").Value = 118
$ws4.Range("F35").Value = 4978
$ws4.Range("F40").Value = 12952
$ws4.Range("F41").Value = 1353
$ws4.Range("F42").Value = 115
$ws4.Range("F45").Value = 296
$ws4.Range("F46").Value = 410
